$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.809.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.25%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.266.40'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.76%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.484'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.56%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.11%  '

$ws.Range('E11').Value = '  -0.38%  '

$ws.Range('E12').Value = '  -1.98%  '

$ws.Range('E13').Value = '  -0.76%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.617.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.09%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.51%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.264.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.16%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.784'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.73%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.746.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.83'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.09%  '

$ws.Range('E20').Value = '  -0.07%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.39%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.63'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.88%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.43%  '

$ws.Range('E25').Value = '  -0.10%  '

$ws.Range('E26').Value = '  +1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '

$ws.Range('E29').Value = '  -8.43%  '

$ws.Range('E30').Value = '  +2.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.61%  '

$ws.Range('E32').Value = '  +2.05%  '

$ws.Range('E33').Value = '  +0.05%  '

$ws.Range('E34').Value = '  +0.30%  '

$ws.Range('E35').Value = '  -1.71%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.77%  '

$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.99%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.105'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.42%  '

$ws.Range('E39').Value = '  +0.70%  '

$ws.Range('E40').Value = '  +0.02%  '

$ws.Range('E41').Value = '  -1.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.61%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.005.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.13%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.17%  '

$ws.Range('E45').Value = '  +1.25%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.78%  '

$ws.Range('E47').Value = '  -1.29%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '73.08'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.47%  '

$ws.Range('E50').Value = '  -0.68%  '

$ws.Range('E51').Value = '  +0.03%  '
